# Append the new OSTEOPOROSE_PATCH0_6 batch entry as row 6 of the "Batches"
# sheet (dimension grows from A1:E5 to A1:E6), mirroring rows 3/4 which carry
# no explicit cell style.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Batches")

# A6 looks like a date ("2026-01-27"); force it to be stored as text so it
# matches the other "Data" column entries (which are plain strings, not
# date values), then drop the formatting Excel auto-applies so the cell
# ends up with no explicit style, same as its neighbours.
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "2026-01-27"
$ws.Range("A6").ClearFormats()

$ws.Range("B6").Value = "OSTEOPOROSE_PATCH0_6"
$ws.Range("C6").Value = "Viewer: fit sem cortes (vh var + floor scale + translate/scale overflow). PDF: novo player full-screen + teclado. Print: print-fit mais robusto."
$ws.Range("D6").Value = "OSTEOPOROSE_PATCH0_6.zip"
$ws.Range("E6").Value = "PDF em assets é placeholder; regenerar via print.html para 72 páginas 16:9."
